$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save"), styled like the other header cells (B1:G1).
# Copy G1's formatting (style index) onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2, corresponding numeric value for the "Save" column.
$ws.Range("H2").Value = 0

$excel.CutCopyMode = 0
